$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: observation about "Manter estoque" wording -------------
$r15 = $ws.Range("B15")
$r15.Value = 'O termo "Manter" para o caso de Estoque não é adequado, pois o CRUD para esse tipo de entidade é bem diferenciado. Assim, ao invés de cadastrar, existe o "Adicionar" ou "Dar Entrada" no estoque. Ao invés de remover, existe o "Dar Baixa" no estoque. Portanto, devem ser casos de uso diferentes.'
$r15.HorizontalAlignment = -4131
$r15.WrapText = $true
$ws.Range("B15:E15").Merge() | Out-Null
$ws.Rows.Item(15).RowHeight = 57.75

# --- New row 16: question about payments ---------------------------------
$ws.Range("B16").Value = 'A parte de pagamentos está fora do sistema?'

# --- New row 17: observation about "manter orçamentos" wording -----------
$r17 = $ws.Range("B17")
$r17.Value = 'Quando vocês descrevem "Relatório de Orçamentos", a idéia que passa é de uma coisa fixa, um relatório. Já o nome "manter orçamentos" dá uma idéia de CRUD. Isso merece uma reflexão e melhoria.'
$r17.HorizontalAlignment = -4131
$r17.VerticalAlignment = -4160
$r17.WrapText = $true
$ws.Range("B17:E17").Merge() | Out-Null
$ws.Rows.Item(17).RowHeight = 48

# --- Selection, matching the saved cursor position in the workbook -------
$null = $ws.Range("B15:E15").Select()
